# Hortaliza, Vega Central Mapocho de Santiago - Orégano
# Insert two new weekly price records into the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 13 (pushes old rows 13-24 down to 14-25) ---
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value  = 9
$ws.Cells.Item(13, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(13, 3).Value  = "Metropolitana"
$ws.Cells.Item(13, 4).Value  = 44426
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 5).Value  = 13
$ws.Cells.Item(13, 6).Value  = 100112029
$ws.Cells.Item(13, 7).Value  = "Orégano"
$ws.Cells.Item(13, 8).Value  = "Sin especificar"
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 16
$ws.Cells.Item(13, 11).Value = 10000
$ws.Cells.Item(13, 12).Value = 10500
$ws.Cells.Item(13, 13).Value = 10250
$ws.Cells.Item(13, 14).Value = "$/docena de atados"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 3417
$ws.Cells.Item(13, 17).Value = 3
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# --- Insert new record at row 20 (pushes rows 20-25 down to 21-26) ---
$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value  = 9
$ws.Cells.Item(20, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value  = "Metropolitana"
$ws.Cells.Item(20, 4).Value  = 44419
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
$ws.Cells.Item(20, 5).Value  = 13
$ws.Cells.Item(20, 6).Value  = 100112029
$ws.Cells.Item(20, 7).Value  = "Orégano"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 16
$ws.Cells.Item(20, 11).Value = 10000
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 13).Value = 10000
$ws.Cells.Item(20, 14).Value = "$/docena de atados"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 3333
$ws.Cells.Item(20, 17).Value = 3
$ws.Cells.Item(20, 18).Value = "Hortaliza"
